$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '58.188.03'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.52%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.475.54'
$ws.Range('D3').ClearFormats()

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '520.31'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.10%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '131.63'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.94%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('E8').Value = '  -1.40%  '

$ws.Range('E9').Value = '  -1.46%  '

$ws.Range('E10').Value = '  -0.71%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.35'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.30%  '

$ws.Range('E12').Value = '  -1.18%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.915.54'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.80%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '58.127.25'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.50%  '

$ws.Range('E15').Value = '  -4.05%  '

$ws.Range('E16').Value = '  -1.72%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.477.37'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.79%  '

$ws.Range('E18').Value = '  -2.61%  '

$ws.Range('E19').Value = '  -2.28%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '321.02'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.00%  '

$ws.Range('E21').Value = '  -0.15%  '

$ws.Range('E22').Value = '  -2.60%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '64.12'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.06%  '

$ws.Range('E24').Value = '  -2.41%  '

$ws.Range('E25').Value = '  -0.09%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.161'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.18%  '

$ws.Range('E27').Value = '  -2.67%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0₃0754'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.81%  '

$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.20'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.23%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.71'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.93%  '

$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.35'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.91%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '166.38'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.39%  '

$ws.Range('E33').Value = '  -0.03%  '

$ws.Range('E34').Value = '  -0.13%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '18.16'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.44%  '

$ws.Range('E36').Value = '  -10.45%  '

$ws.Range('E37').Value = '  -2.44%  '

$ws.Range('E38').Value = '  -3.45%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.795'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.32%  '

$ws.Range('E40').Value = '  -3.71%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '276.62'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.09%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.05'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.09%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.596'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.33%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '126.44'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.76%  '

$ws.Range('E45').Value = '  -2.08%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0491'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.27%  '

$ws.Range('E47').Value = '  -2.47%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '17.21'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.33%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.744.71'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.30%  '

$ws.Range('E51').Value = '  -1.25%  '
